$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("cjf-sa@aomya.cn", "aomya.cn"),
    @("info@casadora.me", "casadora.me"),
    @("info@salt-lamps.com.pk", "salt-lamps.com.pk"),
    @("madsiad@madsiad.org.tr", "madsiad.org.tr"),
    @("marketing@naturalfibres.in", "naturalfibres.in"),
    @("muhasebe@tatmakarna.com", "tatmakarna.com"),
    @("selahattin@korfezreduktor.com", "korfezreduktor.com")
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}
